$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2267135318314786
$ws.Range("D2").Value = 0.02982982747061413
$ws.Range("E2").Value = 0.1342373667421981
$ws.Range("F2").Value = 0.6665497032795145
$ws.Range("G2").Value = 0.5064781280667319
$ws.Range("H2").Value = 0.6495275335726163
$ws.Range("I2").Value = 0.6973097984942953
$ws.Range("K2").Value = 0.8117677129332037
$ws.Range("L2").Value = 0.1567507389160845
$ws.Range("M2").Value = 0.2909352412059008
$ws.Range("N2").Value = 1.20474878347494
$ws.Range("O2").Value = 2.276973025508013
$ws.Range("C3").Value = 0.2242952465281292
$ws.Range("D3").Value = 0.02754505269187035
$ws.Range("E3").Value = 0.1340177325350105
$ws.Range("F3").Value = 0.6654787135903746
$ws.Range("G3").Value = 0.5069613789251832
$ws.Range("H3").Value = 0.6535056646909467
$ws.Range("I3").Value = 0.6986251025391113
$ws.Range("K3").Value = 0.7147362434391198
$ws.Range("L3").Value = 0.1573582303633678
$ws.Range("M3").Value = 0.2690871059606579
$ws.Range("N3").Value = 1.205771604953711
$ws.Range("O3").Value = 2.286060587160748
$ws.Range("C4").Value = 0.2229078568501279
$ws.Range("D4").Value = 0.02613049764488551
$ws.Range("E4").Value = 0.1339472590345601
$ws.Range("F4").Value = 0.6652512479300441
$ws.Range("G4").Value = 0.5076249767163858
$ws.Range("H4").Value = 0.6562470512049075
$ws.Range("I4").Value = 0.6998401724891643
$ws.Range("K4").Value = 0.655016637315839
$ws.Range("L4").Value = 0.1578110089507092
$ws.Range("M4").Value = 0.2557212878434214
$ws.Range("N4").Value = 1.2068107608657
$ws.Range("O4").Value = 2.293034269777735
$ws.Range("C5").Value = 0.2223670729889164
$ws.Range("D5").Value = 0.02555114554907334
$ws.Range("E5").Value = 0.1339347532359731
$ws.Range("F5").Value = 0.6652667180200211
$ws.Range("G5").Value = 0.5079875740612962
$ws.Range("H5").Value = 0.6574393683732964
$ws.Range("I5").Value = 0.7004378056695764
$ws.Range("K5").Value = 0.6306465683619251
$ws.Range("L5").Value = 0.1580155882884995
$ws.Range("M5").Value = 0.2502873459237378
$ws.Range("N5").Value = 1.20733790654721
$ws.Range("O5").Value = 2.296226524939073
$ws.Range("C6").Value = 0.2222787638300332
$ws.Range("D6").Value = 0.0254547697683023
$ws.Range("E6").Value = 0.1339336565471356
$ws.Range("F6").Value = 0.6652758198263129
$ws.Range("G6").Value = 0.5080533479470262
$ws.Range("H6").Value = 0.6576418937307835
$ws.Range("I6").Value = 0.7005432325076377
$ws.Range("K6").Value = 0.6265979469023648
$ws.Range("L6").Value = 0.1580507708477654
$ws.Range("M6").Value = 0.2493858253825465
$ws.Range("N6").Value = 1.207431709594893
$ws.Range("O6").Value = 2.296777758805291
$ws.Range("C7").Value = 0.2229004639787462
$ws.Range("D7").Value = 0.02612269603888251
$ws.Range("E7").Value = 0.133947024701417
$ws.Range("F7").Value = 0.6652510186056375
$ws.Range("G7").Value = 0.5076294937195698
$ws.Range("H7").Value = 0.6562628267519202
$ws.Range("I7").Value = 0.699847817404418
$ws.Range("K7").Value = 0.6546881086454448
$ws.Range("L7").Value = 0.1578136867124691
$ws.Range("M7").Value = 0.2556479516855532
$ws.Range("N7").Value = 1.206817449938569
$ws.Range("O7").Value = 2.293075902906821
$ws.Range("C8").Value = 0.2258595315104799
$ws.Range("D8").Value = 0.02904448229643464
$ws.Range("E8").Value = 0.1341482877197642
$ws.Range("F8").Value = 0.6660911543292869
$ws.Range("G8").Value = 0.5065685581702084
$ws.Range("H8").Value = 0.6508372041210393
$ws.Range("I8").Value = 0.6976787565125306
$ws.Range("K8").Value = 0.7783417826204868
$ws.Range("L8").Value = 0.1569436489400893
$ws.Range("M8").Value = 0.28339206767712
$ws.Range("N8").Value = 1.205016294942979
$ws.Range("O8").Value = 2.279817064209311
$ws.Range("C9").Value = 0.2324322410625257
$ws.Range("D9").Value = 0.03468018548072394
$ws.Range("E9").Value = 0.135052971163578
$ws.Range("F9").Value = 0.6711524178540813
$ws.Range("G9").Value = 0.507403388753346
$ws.Range("H9").Value = 0.6425666678160269
$ws.Range("I9").Value = 0.6966581842316089
$ws.Range("K9").Value = 1.01963119509486
$ws.Range("L9").Value = 0.1558702632572775
$ws.Range("M9").Value = 0.338171934971335
$ws.Range("N9").Value = 1.204734468291235
$ws.Range("O9").Value = 2.264882124864442
$ws.Range("C10").Value = 0.2377270509038141
$ws.Range("D10").Value = 0.03876233709664234
$ws.Range("E10").Value = 0.1360276850423503
$ws.Range("F10").Value = 0.6769549459452264
$ws.Range("G10").Value = 0.509800992707369
$ws.Range("H10").Value = 0.6379325791645698
$ws.Range("I10").Value = 0.6978802863579219
$ws.Range("K10").Value = 1.196101738968252
$ws.Range("L10").Value = 0.1554672489981073
$ws.Range("M10").Value = 0.3786300890513559
$ws.Range("N10").Value = 1.206494310193463
$ws.Range("O10").Value = 2.260666387310721
$ws.Range("C11").Value = 0.2402362391518835
$ws.Range("D11").Value = 0.04060650903184637
$ws.Range("E11").Value = 0.1365382620066669
$ws.Range("F11").Value = 0.6800479047081254
$ws.Range("G11").Value = 0.511280715595916
$ws.Range("H11").Value = 0.6361371980709407
$ws.Range("I11").Value = 0.6988647284059155
$ws.Range("K11").Value = 1.276192870376519
$ws.Range("L11").Value = 0.1553676191265581
$ws.Range("M11").Value = 0.397078205163389
$ws.Range("N11").Value = 1.207719107810505
$ws.Range("O11").Value = 2.260218319144741
$ws.Range("C12").Value = 0.2412007871752877
$ws.Range("D12").Value = 0.0413029780504246
$ws.Range("E12").Value = 0.1367412445896967
$ws.Range("F12").Value = 0.6812843437103311
$ws.Range("G12").Value = 0.5118970957497595
$ws.Range("H12").Value = 0.6355022631312437
$ws.Range("I12").Value = 0.6992991303289386
$ws.Range("K12").Value = 1.306492873477225
$ws.Range("L12").Value = 0.1553419239698357
$ws.Range("M12").Value = 0.4040699178447511
$ws.Range("N12").Value = 1.208243656920885
$ws.Range("O12").Value = 2.260260140754156
$ws.Range("C13").Value = 0.2409924166403528
$ws.Range("D13").Value = 0.04115306511494765
$ws.Range("E13").Value = 0.1366971003416566
$ws.Range("F13").Value = 0.6810151543993541
$ws.Range("G13").Value = 0.5117618532681263
$ws.Range("H13").Value = 0.6356370095861621
$ws.Range("I13").Value = 0.6992028338117038
$ws.Range("K13").Value = 1.29996853694314
$ws.Range("L13").Value = 0.1553469227960207
$ws.Range("M13").Value = 0.4025638760128203
$ws.Range("N13").Value = 1.208127988326268
$ws.Range("O13").Value = 2.260241724971195
$ws.Range("C14").Value = 0.2403153055391982
$ws.Range("D14").Value = 0.04066384581604865
$ws.Range("E14").Value = 0.1365547684660839
$ws.Range("F14").Value = 0.6801483204785868
$ws.Range("G14").Value = 0.5113303018588198
$ws.Range("H14").Value = 0.6360840611458798
$ws.Range("I14").Value = 0.6988992320756182
$ws.Range("K14").Value = 1.278686257873119
$ws.Range("L14").Value = 0.155365264044093
$ws.Range("M14").Value = 0.3976533036423859
$ws.Range("N14").Value = 1.207761046835159
$ws.Range("O14").Value = 2.260217520273955
$ws.Range("C15").Value = 0.2399024246465729
$ws.Range("D15").Value = 0.04036393913210645
$ws.Range("E15").Value = 0.1364688406172725
$ws.Range("F15").Value = 0.6796258512547979
$ws.Range("G15").Value = 0.5110732652522501
$ws.Range("H15").Value = 0.6363637444337371
$ws.Range("I15").Value = 0.6987212912457537
$ws.Range("K15").Value = 1.265646441587762
$ws.Range("L15").Value = 0.1553780654581942
$ws.Range("M15").Value = 0.3946461799434857
$ws.Range("N15").Value = 1.207544187720643
$ws.Range("O15").Value = 2.260230241687282
$ws.Range("C16").Value = 0.2375650855395151
$ws.Range("D16").Value = 0.03864155513504386
$ws.Range("E16").Value = 0.1359956675307998
$ws.Range("F16").Value = 0.6767619374678517
$ws.Range("G16").Value = 0.5097121262093225
$ws.Range("H16").Value = 0.638056201304309
$ws.Range("I16").Value = 0.6978245730744703
$ws.Range("K16").Value = 1.190863677704613
$ws.Range("L16").Value = 0.1554754440362274
$ws.Range("M16").Value = 0.3774252986572222
$ws.Range("N16").Value = 1.206422781183647
$ws.Range("O16").Value = 2.260725256933114
$ws.Range("C17").Value = 0.2361568961797076
$ws.Range("D17").Value = 0.0375816197868204
$ws.Range("E17").Value = 0.135722580978797
$ws.Range("F17").Value = 0.6751211430593216
$ws.Range("G17").Value = 0.5089768240835468
$ws.Range("H17").Value = 0.6391745364681327
$ws.Range("I17").Value = 0.6973842132840034
$ws.Range("K17").Value = 1.144937762706434
$ws.Range("L17").Value = 0.1555566183108752
$ws.Range("M17").Value = 0.3668716711290259
$ws.Range("N17").Value = 1.205843305263926
$ws.Range("O17").Value = 2.261405467671523
$ws.Range("C18").Value = 0.2353564125794207
$ws.Range("D18").Value = 0.03697076851771897
$ws.Range("E18").Value = 0.1355718329114524
$ws.Range("F18").Value = 0.6742200733602672
$ws.Range("G18").Value = 0.5085905103653516
$ws.Range("H18").Value = 0.639847206688998
$ws.Range("I18").Value = 0.6971712628167808
$ws.Range("K18").Value = 1.118504983624803
$ws.Range("L18").Value = 0.155611186401778
$ws.Range("M18").Value = 0.3608056245216886
$ws.Range("N18").Value = 1.205549952416533
$ws.Range("O18").Value = 2.261935032164558
$ws.Range("C19").Value = 0.2350870116441399
$ws.Range("D19").Value = 0.03676373898171192
$ws.Range("E19").Value = 0.1355218790281576
$ws.Range("F19").Value = 0.6739223152420877
$ws.Range("G19").Value = 0.5084659964335998
$ws.Range("H19").Value = 0.6400800173199741
$ws.Range("I19").Value = 0.6971060887261515
$ws.Range("K19").Value = 1.109552369889684
$ws.Range("L19").Value = 0.1556310155285558
$ws.Range("M19").Value = 0.3587524866465586
$ws.Range("N19").Value = 1.205457497014805
$ws.Range("O19").Value = 2.262138086318657
$ws.Range("C20").Value = 0.2363058208912463
$ws.Range("D20").Value = 0.03769457663173625
$ws.Range("E20").Value = 0.135750997143802
$ws.Range("F20").Value = 0.6752913919711148
$ws.Range("G20").Value = 0.5090513082701307
$ws.Range("H20").Value = 0.6390524417303567
$ws.Range("I20").Value = 0.6974269160562372
$ws.Range("K20").Value = 1.149828470824616
$ws.Range("L20").Value = 0.1555471617733986
$ws.Range("M20").Value = 0.3679946995215815
$ws.Range("N20").Value = 1.20590085906305
$ws.Range("O20").Value = 2.261318740794309
$ws.Range("C21").Value = 0.240513800096906
$ws.Range("D21").Value = 0.04080759265809064
$ws.Range("E21").Value = 0.1365963133826114
$ws.Range("F21").Value = 0.6804011609597325
$ws.Range("G21").Value = 0.5114555372554577
$ws.Range("H21").Value = 0.6359515319092424
$ws.Range("I21").Value = 0.6989867350772556
$ws.Range("K21").Value = 1.284938175813068
$ws.Range("L21").Value = 0.1553595502478302
$ws.Range("M21").Value = 0.3990955035929886
$ws.Range("N21").Value = 1.207867179894578
$ws.Range("O21").Value = 2.260218888635904
$ws.Range("C22").Value = 0.2433476937492003
$ws.Range("D22").Value = 0.0428311606013807
$ws.Range("E22").Value = 0.1372049449892856
$ws.Range("F22").Value = 0.6841207448737237
$ws.Range("G22").Value = 0.5133535352297258
$ws.Range("H22").Value = 0.6341868220801246
$ws.Range("I22").Value = 0.7003653048190941
$ws.Range("K22").Value = 1.373071561342101
$ws.Range("L22").Value = 0.1553070680745634
$ws.Range("M22").Value = 0.4194553451254066
$ws.Range("N22").Value = 1.209506234567058
$ws.Range("O22").Value = 2.260732898011184
$ws.Range("C23").Value = 0.2418275581263316
$ws.Range("D23").Value = 0.04175215960966483
$ws.Range("E23").Value = 0.1368749745397935
$ws.Range("F23").Value = 0.6821007549937903
$ws.Range("G23").Value = 0.5123106134607696
$ws.Range("H23").Value = 0.6351047248399624
$ws.Range("I23").Value = 0.6995966762781762
$ws.Range("K23").Value = 1.326049201163471
$ws.Range("L23").Value = 0.1553286629637185
$ws.Range("M23").Value = 0.4085859803579837
$ws.Range("N23").Value = 1.208599138325241
$ws.Range("O23").Value = 2.260345707315736
$ws.Range("C24").Value = 0.2362384636685135
$ws.Range("D24").Value = 0.03764351345101602
$ws.Range("E24").Value = 0.1357381307167032
$ws.Range("F24").Value = 0.6752142908497376
$ws.Range("G24").Value = 0.5090175205275216
$ws.Range("H24").Value = 0.6391075481691217
$ws.Range("I24").Value = 0.6974074848561003
$ws.Range("K24").Value = 1.147617472981551
$ws.Range("L24").Value = 0.1555514124665649
$ws.Range("M24").Value = 0.3674869740705446
$ws.Range("N24").Value = 1.205874715009841
$ws.Range("O24").Value = 2.261357518602154
$ws.Range("C25").Value = 0.2305720847177639
$ws.Range("D25").Value = 0.03316575085943185
$ws.Range("E25").Value = 0.1347537203065876
$ws.Range("F25").Value = 0.6694174410698395
$ws.Range("G25").Value = 0.5068647121508576
$ws.Range("H25").Value = 0.6445506451218961
$ws.Range("I25").Value = 0.6965880580248296
$ws.Range("K25").Value = 0.9544918535481486
$ws.Range("L25").Value = 0.1560929161708771
$ws.Range("M25").Value = 0.323314220886445
$ws.Range("N25").Value = 1.204464272424644
$ws.Range("O25").Value = 2.267736619077425
